$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (I0, IF) - copy the header style (bold/border/center) from H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data columns I (I0) and J (IF) for rows 2-16
$iValues = @(1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 4, 1)
$jValues = @(2, 3, 6, 5, 6, 4, 6, 6, 5, 6, 5, 5, 3, 6, 2)

for ($k = 0; $k -lt 15; $k++) {
    $row = $k + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$k]
    $ws.Cells.Item($row, 10).Value = $jValues[$k]
}
